# Update scripts with new TPM: recalculated ligand/receptor expression and
# derived specificity/weight columns (G,H,I,J,M,N,O,P,Q,R,S,T) for every
# sending/target cluster pair row in the Adam12-Itgb1 LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.2162966666666667
$ws.Range("H2").Value = 0.6488900000000001
$ws.Range("I2").Value = 0.02888548604596741
$ws.Range("J2").Value = 0.0288854860459674
$ws.Range("M2").Value = 61.04160633333334
$ws.Range("N2").Value = 183.124819
$ws.Range("O2").Value = 0.2043613460574534
$ws.Range("P2").Value = 0.2043613460574534
$ws.Range("Q2").Value = 13.20309597787889
$ws.Range("R2").Value = 118.82786380091
$ws.Range("S2").Value = 0.005903076809877687
$ws.Range("T2").Value = 0.005903076809877687

$ws.Range("G3").Value = 0.2162966666666667
$ws.Range("H3").Value = 0.6488900000000001
$ws.Range("I3").Value = 0.02888548604596741
$ws.Range("J3").Value = 0.0288854860459674
$ws.Range("O3").Value = 0.3559304658284363
$ws.Range("P3").Value = 0.3559304658284363
$ws.Range("Q3").Value = 22.99546461424667
$ws.Range("R3").Value = 206.95918152822
$ws.Range("S3").Value = 0.01028122450402198
$ws.Range("T3").Value = 0.01028122450402198

$ws.Range("G4").Value = 0.2162966666666667
$ws.Range("H4").Value = 0.6488900000000001
$ws.Range("I4").Value = 0.02888548604596741
$ws.Range("J4").Value = 0.0288854860459674
$ws.Range("M4").Value = 131.3384093333333
$ws.Range("N4").Value = 394.015228
$ws.Range("O4").Value = 0.4397081881141102
$ws.Range("P4").Value = 0.4397081881141103
$ws.Range("Q4").Value = 28.40806014410222
$ws.Range("R4").Value = 255.67254129692
$ws.Range("S4").Value = 0.01270118473206774
$ws.Range("T4").Value = 0.01270118473206774

$ws.Range("I5").Value = 0.7789723686414617
$ws.Range("J5").Value = 0.7789723686414615
$ws.Range("M5").Value = 61.04160633333334
$ws.Range("N5").Value = 183.124819
$ws.Range("O5").Value = 0.2043613460574534
$ws.Range("P5").Value = 0.2043613460574534
$ws.Range("Q5").Value = 356.0558728671524
$ws.Range("R5").Value = 3204.502855804372
$ws.Range("S5").Value = 0.1591918417971319
$ws.Range("T5").Value = 0.1591918417971319

$ws.Range("I6").Value = 0.7789723686414617
$ws.Range("J6").Value = 0.7789723686414615
$ws.Range("O6").Value = 0.3559304658284363
$ws.Range("P6").Value = 0.3559304658284363
$ws.Range("S6").Value = 0.2772599980380359
$ws.Range("T6").Value = 0.2772599980380358

$ws.Range("I7").Value = 0.7789723686414617
$ws.Range("J7").Value = 0.7789723686414615
$ws.Range("M7").Value = 131.3384093333333
$ws.Range("N7").Value = 394.015228
$ws.Range("O7").Value = 0.4397081881141102
$ws.Range("P7").Value = 0.4397081881141103
$ws.Range("Q7").Value = 766.0973356565613
$ws.Range("R7").Value = 6894.876020909052
$ws.Range("S7").Value = 0.3425205288062939
$ws.Range("T7").Value = 0.3425205288062938

$ws.Range("G8").Value = 1.438774666666667
$ws.Range("H8").Value = 4.316324
$ws.Range("I8").Value = 0.192142145312571
$ws.Range("J8").Value = 0.192142145312571
$ws.Range("M8").Value = 61.04160633333334
$ws.Range("N8").Value = 183.124819
$ws.Range("O8").Value = 0.2043613460574534
$ws.Range("P8").Value = 0.2043613460574534
$ws.Range("Q8").Value = 87.82511680503956
$ws.Range("R8").Value = 790.4260512453561
$ws.Range("S8").Value = 0.03926642745044382
$ws.Range("T8").Value = 0.03926642745044382

$ws.Range("G9").Value = 1.438774666666667
$ws.Range("H9").Value = 4.316324
$ws.Range("I9").Value = 0.192142145312571
$ws.Range("J9").Value = 0.192142145312571
$ws.Range("O9").Value = 0.3559304658284363
$ws.Range("P9").Value = 0.3559304658284363
$ws.Range("Q9").Value = 152.9625603809947
$ws.Range("R9").Value = 1376.663043428952
$ws.Range("S9").Value = 0.06838924328637849
$ws.Range("T9").Value = 0.06838924328637849

$ws.Range("G10").Value = 1.438774666666667
$ws.Range("H10").Value = 4.316324
$ws.Range("I10").Value = 0.192142145312571
$ws.Range("J10").Value = 0.192142145312571
$ws.Range("M10").Value = 131.3384093333333
$ws.Range("N10").Value = 394.015228
$ws.Range("O10").Value = 0.4397081881141102
$ws.Range("P10").Value = 0.4397081881141103
$ws.Range("Q10").Value = 188.9663761090968
$ws.Range("R10").Value = 1700.697384981872
$ws.Range("S10").Value = 0.08448647457574866
$ws.Range("T10").Value = 0.08448647457574868
